# Refresh the crypto price/volume table to the latest scrape.
# Note: several "Price" cells look like plain numbers (e.g. "210.06") but must
# stay text, matching the rest of the column (which uses dotted "thousand"
# groups like "26.291.26" that can't be numbers anyway). A leading apostrophe
# forces Excel to store them as text instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.291.26"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "1.591.07"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").Value = "'210.06"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "'0.0844"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "1.814.55"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.07"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.577.66"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "'64.63"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").Value = "26.306.32"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("D20").Value = "'211.61"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").Value = "'4.27"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "'145.44"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").Value = "'15.28"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "'3.22"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "1.300.42"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("D35").Value = "'0.617"
$ws.Range("E35").Value = "  +4.21%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  -13.88%  "
$ws.Range("D40").Value = "'0.813"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("D44").Value = "'62.64"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "1.727.83"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").Value = "'88.58"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("D50").Value = "'0.0986"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("D51").Value = "'0.0505"
$ws.Range("E51").Value = "  -1.33%  "
